# "changes on the dataset"
#
# The author re-saved this workbook from a different machine/Excel build
# (the diff is full of Excel-version fingerprint churn: fileVersion,
# xr:revisionPtr, bookViews window coords, defaultThemeVersion, x14ac
# dyDescent stamps, best-fit column-width epsilons, etc. — none of which
# are meaningful, user-driven edits and none of which are reachable
# through the object model). The two real content changes are:
#
#   1. The worksheet was renamed from
#      "Template emotion_datasheet1_use" to "Sheet1".
#   2. A hidden workbook-scoped defined name, LOCAL_MYSQL_DATE_FORMAT
#      (a helper formula that ships with the "MySQL for Excel" add-in),
#      was (re)created.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the sheet.
$ws.Name = "Sheet1"

# 2) Recreate the hidden defined name with its MySQL-for-Excel formula.
$formula = "=REPT(LOCAL_YEAR_FORMAT,4)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_MONTH_FORMAT,2)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_DAY_FORMAT,2)&"" ""&REPT(LOCAL_HOUR_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_MINUTE_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_SECOND_FORMAT,2)"
$definedName = $wb.Names.Add("LOCAL_MYSQL_DATE_FORMAT", $formula)
$definedName.Visible = $false
